# "Added year to archetypes (only one construction per function for now)
#  - correction of Temperature set points to low values for no-heating
#  system case and high values for no-cooling system case (input has to be
#  double). I.e., can not be 'None' or 'NaN' or similar"
#
#  - ARCHITECTURE and HVAC sheets: fill the (previously blank) B/C columns
#    (year_start / year_end) for every data row (2-19) with 1950 / 2030.
#  - INDOOR_COMFORT sheet: replace the placeholder text "None" used for the
#    setback temperatures (columns C, E; and for rows 13 & 15 also B, D)
#    with real numeric values (10 for the setback columns, 50 for the B/D
#    setpoint columns on rows 13 & 15) since downstream code expects a
#    double and chokes on 'None'/'NaN'.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ARCHITECTURE + HVAC: set year_start (B) / year_end (C) for rows 2..19
# ---------------------------------------------------------------------------
$wsArchitecture = $wb.Worksheets.Item("ARCHITECTURE")
$wsHvac = $wb.Worksheets.Item("HVAC")

foreach ($ws in @($wsArchitecture, $wsHvac)) {
    for ($row = 2; $row -le 19; $row++) {
        $ws.Cells.Item($row, 2).Value = 1950
        $ws.Cells.Item($row, 3).Value = 2030
    }
}

# ---------------------------------------------------------------------------
# INDOOR_COMFORT: replace 'None' placeholders with real numbers
# ---------------------------------------------------------------------------
$wsComfort = $wb.Worksheets.Item("INDOOR_COMFORT")

# Columns C (Tcs_setb_C) and E (Ths_setb_C) are 'None' on every data row -
# correct them all to 10.
for ($row = 2; $row -le 19; $row++) {
    $wsComfort.Cells.Item($row, 3).Value = 10
    $wsComfort.Cells.Item($row, 5).Value = 10
}

# Rows 13 (SCHOOL) and 15 (GYM) also had 'None' in B (Ths_set_C) and
# D (Tcs_set_C) - set those to 50.
foreach ($row in @(13, 15)) {
    $wsComfort.Cells.Item($row, 2).Value = 50
    $wsComfort.Cells.Item($row, 4).Value = 50
}

# ---------------------------------------------------------------------------
# Restore the view/selection state recorded in the edit, without disturbing
# which sheet/tab ends up active (ARCHITECTURE stays the active tab, as in
# the original workbook).
# ---------------------------------------------------------------------------
$wsHvac.Range("B19:C19").Select()
$wsComfort.Range("H36").Select()
$wsArchitecture.Range("B19:C19").Select()
